# Add a new "password_notify_header" mail template row to the
# mail_template sheet (row 13), reusing the existing "Password
# notification" (column C) and "body" (column D) shared strings, and
# introducing three brand-new shared strings for the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "password_notify_header"
$ws.Range("C13").Value = "Password notification"
$ws.Range("D13").Value = "body"
$ws.Range("E13").Value = '[${system:site_name}]Password notification header'
$ws.Range("F13").Value = '***************************************************************************\nThis email is automatically encrypted as an attachment.\nYou will receive a password from the sender later.\n***************************************************************************\n\n'

# Match the author's final selection, which lands the cursor on the new
# row's last populated cell.
$ws.Range("F13").Select()
